$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.019999999999999
$ws.Cells.Item(2, 3).Value = 1.032986268198313
$ws.Cells.Item(2, 4).Value = 1.035166563004088
$ws.Cells.Item(2, 5).Value = 1.041911371863755
$ws.Cells.Item(2, 6).Value = 1.051950751342542
$ws.Cells.Item(2, 9).Value = 1.031427523257402
$ws.Cells.Item(2, 10).Value = 1.038113475487467
$ws.Cells.Item(2, 11).Value = 1.037963946410103
$ws.Cells.Item(2, 12).Value = 1.044689548032033
$ws.Cells.Item(2, 13).Value = 1.054700832904487
$ws.Cells.Item(2, 14).Value = 1.01658377132681
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034249849746943
$ws.Cells.Item(3, 4).Value = 1.036098652130245
$ws.Cells.Item(3, 5).Value = 1.043108013993592
$ws.Cells.Item(3, 6).Value = 1.05340289474558
$ws.Cells.Item(3, 9).Value = 1.031643089641752
$ws.Cells.Item(3, 10).Value = 1.039017988224845
$ws.Cells.Item(3, 11).Value = 1.038705047700663
$ws.Cells.Item(3, 12).Value = 1.045695902418718
$ws.Cells.Item(3, 13).Value = 1.055964082948603
$ws.Cells.Item(3, 14).Value = 1.016891237323633
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.035066822276832
$ws.Cells.Item(4, 4).Value = 1.036700926404075
$ws.Cells.Item(4, 5).Value = 1.043882093511119
$ws.Cells.Item(4, 6).Value = 1.054342707907252
$ws.Cells.Item(4, 9).Value = 1.031780855164349
$ws.Cells.Item(4, 10).Value = 1.0396021543089
$ws.Cells.Item(4, 11).Value = 1.039183119822329
$ws.Cells.Item(4, 12).Value = 1.046346287273835
$ws.Cells.Item(4, 13).Value = 1.056781153963136
$ws.Cells.Item(4, 14).Value = 1.017089616823679
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035410125374846
$ws.Cells.Item(5, 4).Value = 1.036953920535579
$ws.Cells.Item(5, 5).Value = 1.044207464577027
$ws.Cells.Item(5, 6).Value = 1.054737851951341
$ws.Cells.Item(5, 9).Value = 1.031838360383018
$ws.Cells.Item(5, 10).Value = 1.039847473061321
$ws.Cells.Item(5, 11).Value = 1.039383750886628
$ws.Cells.Item(5, 12).Value = 1.046619521412014
$ws.Cells.Item(5, 13).Value = 1.057124573293156
$ws.Cells.Item(5, 14).Value = 1.017172879321063
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035467758648192
$ws.Cells.Item(6, 4).Value = 1.036996387561658
$ws.Cells.Item(6, 5).Value = 1.044262092839255
$ws.Cells.Item(6, 6).Value = 1.054804201213736
$ws.Cells.Item(6, 9).Value = 1.03184799165006
$ws.Cells.Item(6, 10).Value = 1.039888647698626
$ws.Cells.Item(6, 11).Value = 1.03941741721646
$ws.Cells.Item(6, 12).Value = 1.046665387690506
$ws.Cells.Item(6, 13).Value = 1.05718223046745
$ws.Cells.Item(6, 14).Value = 1.017186851488328
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.035071410103421
$ws.Cells.Item(7, 4).Value = 1.036704307715875
$ws.Cells.Item(7, 5).Value = 1.043886441336164
$ws.Cells.Item(7, 6).Value = 1.054347987653668
$ws.Cells.Item(7, 9).Value = 1.031781625167898
$ws.Cells.Item(7, 10).Value = 1.0396054333075
$ws.Cells.Item(7, 11).Value = 1.039185802039688
$ws.Cells.Item(7, 12).Value = 1.046349938976186
$ws.Cells.Item(7, 13).Value = 1.056785743048924
$ws.Cells.Item(7, 14).Value = 1.017090729915865
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033413437897558
$ws.Cells.Item(8, 4).Value = 1.035481743704682
$ws.Cells.Item(8, 5).Value = 1.042315831448132
$ws.Cells.Item(8, 6).Value = 1.05244147509615
$ws.Cells.Item(8, 9).Value = 1.031500731404686
$ws.Cells.Item(8, 10).Value = 1.038419391733524
$ws.Cells.Item(8, 11).Value = 1.038214710308395
$ws.Cells.Item(8, 12).Value = 1.045029815980832
$ws.Cells.Item(8, 13).Value = 1.055127825702143
$ws.Cells.Item(8, 14).Value = 1.016687799600382
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030486759789462
$ws.Cells.Item(9, 4).Value = 1.033320856116463
$ws.Cells.Item(9, 5).Value = 1.039546352888951
$ws.Cells.Item(9, 6).Value = 1.049083164722806
$ws.Cells.Item(9, 9).Value = 1.030992564562541
$ws.Cells.Item(9, 10).Value = 1.036320809255022
$ws.Cells.Item(9, 11).Value = 1.036492198446624
$ws.Cells.Item(9, 12).Value = 1.042697405241559
$ws.Cells.Item(9, 13).Value = 1.052203645993665
$ws.Cells.Item(9, 14).Value = 1.0159733808226
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028531991541223
$ws.Cells.Item(10, 4).Value = 1.031875749828289
$ws.Cells.Item(10, 5).Value = 1.03769861467391
$ws.Cells.Item(10, 6).Value = 1.046844871277783
$ws.Cells.Item(10, 9).Value = 1.030644884502886
$ws.Cells.Item(10, 10).Value = 1.034915823005265
$ws.Cells.Item(10, 11).Value = 1.035336148870178
$ws.Cells.Item(10, 12).Value = 1.041138155255976
$ws.Cells.Item(10, 13).Value = 1.050252164685198
$ws.Cells.Item(10, 14).Value = 1.015494102358212
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02768464283048
$ws.Cells.Item(11, 4).Value = 1.031248910058366
$ws.Cells.Item(11, 5).Value = 1.036898149275331
$ws.Cells.Item(11, 6).Value = 1.04587575080701
$ws.Cells.Item(11, 9).Value = 1.030492216656582
$ws.Cells.Item(11, 10).Value = 1.034306011602713
$ws.Cells.Item(11, 11).Value = 1.034833716521374
$ws.Cells.Item(11, 12).Value = 1.040461926432175
$ws.Cells.Item(11, 13).Value = 1.04940662531169
$ws.Cells.Item(11, 14).Value = 1.015285849282287
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.0273697572679
$ws.Cells.Item(12, 4).Value = 1.031015906710282
$ws.Cells.Item(12, 5).Value = 1.036600759969368
$ws.Cells.Item(12, 6).Value = 1.045515783050175
$ws.Cells.Item(12, 9).Value = 1.030435189827482
$ws.Cells.Item(12, 10).Value = 1.034079281008213
$ws.Cells.Item(12, 11).Value = 1.034646810066543
$ws.Cells.Item(12, 12).Value = 1.040210582155036
$ws.Cells.Item(12, 13).Value = 1.0490924699562
$ws.Cells.Item(12, 14).Value = 1.015208385493123
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027437307805962
$ws.Cells.Item(13, 4).Value = 1.03106589431551
$ws.Cells.Item(13, 5).Value = 1.036664553790902
$ws.Cells.Item(13, 6).Value = 1.045592997124983
$ws.Cells.Item(13, 9).Value = 1.030447436730971
$ws.Cells.Item(13, 10).Value = 1.034127925465872
$ws.Cells.Item(13, 11).Value = 1.034686914870197
$ws.Cells.Item(13, 12).Value = 1.040264503770864
$ws.Cells.Item(13, 13).Value = 1.049159861229268
$ws.Cells.Item(13, 14).Value = 1.015225006698938
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027658617223031
$ws.Cells.Item(14, 4).Value = 1.031229653350332
$ws.Cells.Item(14, 5).Value = 1.036873568242548
$ws.Cells.Item(14, 6).Value = 1.04584599564912
$ws.Cells.Item(14, 9).Value = 1.03048750931647
$ws.Cells.Item(14, 10).Value = 1.034287274474085
$ws.Cells.Item(14, 11).Value = 1.0348182725134
$ws.Cells.Item(14, 12).Value = 1.040441153574154
$ws.Cells.Item(14, 13).Value = 1.049380658876586
$ws.Cells.Item(14, 14).Value = 1.015279448336148
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027794954249461
$ws.Cells.Item(15, 4).Value = 1.031330528487545
$ws.Cells.Item(15, 5).Value = 1.037002340785579
$ws.Cells.Item(15, 6).Value = 1.046001877144919
$ws.Cells.Item(15, 9).Value = 1.030512157042455
$ws.Cells.Item(15, 10).Value = 1.034385425492527
$ws.Cells.Item(15, 11).Value = 1.034899169036536
$ws.Cells.Item(15, 12).Value = 1.040549971686759
$ws.Cells.Item(15, 13).Value = 1.049516688292806
$ws.Cells.Item(15, 14).Value = 1.015312977115036
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028588207524672
$ws.Cells.Item(16, 4).Value = 1.031917327789619
$ws.Cells.Item(16, 5).Value = 1.037751730581286
$ws.Cells.Item(16, 6).Value = 1.046909189808487
$ws.Cells.Item(16, 9).Value = 1.030654971839827
$ws.Cells.Item(16, 10).Value = 1.034956263547091
$ws.Cells.Item(16, 11).Value = 1.035369454393803
$ws.Cells.Item(16, 12).Value = 1.041183011674851
$ws.Cells.Item(16, 13).Value = 1.050308268734535
$ws.Cells.Item(16, 14).Value = 1.015507908141915
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029085544816272
$ws.Cells.Item(17, 4).Value = 1.032285115788914
$ws.Cells.Item(17, 5).Value = 1.03822169822583
$ws.Cells.Item(17, 6).Value = 1.047478340015742
$ws.Cells.Item(17, 9).Value = 1.030743987707263
$ws.Cells.Item(17, 10).Value = 1.035313946934432
$ws.Cells.Item(17, 11).Value = 1.035663953876544
$ws.Cells.Item(17, 12).Value = 1.041579814122967
$ws.Cells.Item(17, 13).Value = 1.050804660159568
$ws.Cells.Item(17, 14).Value = 1.015629989270733
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029375544646529
$ws.Cells.Item(18, 4).Value = 1.032499534385756
$ws.Cells.Item(18, 5).Value = 1.038495785717937
$ws.Cells.Item(18, 6).Value = 1.047810323285475
$ws.Cells.Item(18, 9).Value = 1.030795704645442
$ws.Cells.Item(18, 10).Value = 1.035522438548721
$ws.Cells.Item(18, 11).Value = 1.035835551572518
$ws.Cells.Item(18, 12).Value = 1.041811159892211
$ws.Cells.Item(18, 13).Value = 1.051094145522812
$ws.Cells.Item(18, 14).Value = 1.015701127487368
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029474412190899
$ws.Cells.Item(19, 4).Value = 1.032572627680166
$ws.Cells.Item(19, 5).Value = 1.038589236347893
$ws.Cells.Item(19, 6).Value = 1.047923522425471
$ws.Cells.Item(19, 9).Value = 1.030813304116361
$ws.Cells.Item(19, 10).Value = 1.035593505295812
$ws.Cells.Item(19, 11).Value = 1.035894031672245
$ws.Cells.Item(19, 12).Value = 1.041890025545902
$ws.Cells.Item(19, 13).Value = 1.051192843959849
$ws.Cells.Item(19, 14).Value = 1.015725371992653
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02903219442696
$ws.Cells.Item(20, 4).Value = 1.032245666588642
$ws.Cells.Item(20, 5).Value = 1.038171278965073
$ws.Cells.Item(20, 6).Value = 1.047417274835204
$ws.Cells.Item(20, 9).Value = 1.030734458305938
$ws.Cells.Item(20, 10).Value = 1.035275585295484
$ws.Cells.Item(20, 11).Value = 1.035632375385696
$ws.Cells.Item(20, 12).Value = 1.041537251588855
$ws.Cells.Item(20, 13).Value = 1.050751407399767
$ws.Cells.Item(20, 14).Value = 1.015616898330354
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027593451083292
$ws.Cells.Item(21, 4).Value = 1.031181435026421
$ws.Cells.Item(21, 5).Value = 1.03681202040183
$ws.Cells.Item(21, 6).Value = 1.045771493740565
$ws.Cells.Item(21, 9).Value = 1.03047571775557
$ws.Cells.Item(21, 10).Value = 1.034240356243374
$ws.Cells.Item(21, 11).Value = 1.034779598726402
$ws.Cells.Item(21, 12).Value = 1.040389139128025
$ws.Cells.Item(21, 13).Value = 1.049315641834275
$ws.Cells.Item(21, 14).Value = 1.015263419655261
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026688029385771
$ws.Cells.Item(22, 4).Value = 1.030511342336554
$ws.Cells.Item(22, 5).Value = 1.035957047063707
$ws.Cells.Item(22, 6).Value = 1.044736762533317
$ws.Cells.Item(22, 9).Value = 1.030311190209265
$ws.Cells.Item(22, 10).Value = 1.033588195423617
$ws.Cells.Item(22, 11).Value = 1.034241799104545
$ws.Cells.Item(22, 12).Value = 1.039666331929625
$ws.Cells.Item(22, 13).Value = 1.048412428824665
$ws.Cells.Item(22, 14).Value = 1.015040540665039
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027168090298826
$ws.Cells.Item(23, 4).Value = 1.030866663584071
$ws.Cells.Item(23, 5).Value = 1.036410319236707
$ws.Cells.Item(23, 6).Value = 1.045285291181911
$ws.Cells.Item(23, 9).Value = 1.030398584684002
$ws.Cells.Item(23, 10).Value = 1.033934039609305
$ws.Cells.Item(23, 11).Value = 1.034527051563045
$ws.Cells.Item(23, 12).Value = 1.040049596134669
$ws.Cells.Item(23, 13).Value = 1.04889128693413
$ws.Cells.Item(23, 14).Value = 1.01515875329685
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029056301445025
$ws.Cells.Item(24, 4).Value = 1.032263492311452
$ws.Cells.Item(24, 5).Value = 1.038194061371979
$ws.Cells.Item(24, 6).Value = 1.047444867538138
$ws.Cells.Item(24, 9).Value = 1.030738764864118
$ws.Cells.Item(24, 10).Value = 1.035292919698437
$ws.Cells.Item(24, 11).Value = 1.035646644899073
$ws.Cells.Item(24, 12).Value = 1.041556484082184
$ws.Cells.Item(24, 13).Value = 1.050775470188198
$ws.Cells.Item(24, 14).Value = 1.015622813777985
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.031244004210015
$ws.Cells.Item(25, 4).Value = 1.033880286071674
$ws.Cells.Item(25, 5).Value = 1.040262567370642
$ws.Cells.Item(25, 6).Value = 1.049951249118706
$ws.Cells.Item(25, 9).Value = 1.031125504733481
$ws.Cells.Item(25, 10).Value = 1.036864378843645
$ws.Cells.Item(25, 11).Value = 1.0369388607175
$ws.Cells.Item(25, 12).Value = 1.043301137801981
$ws.Cells.Item(25, 13).Value = 1.052959959230406
$ws.Cells.Item(25, 14).Value = 1.016158601023321
